# Replace the placeholder constructor names in column A with the canonical
# names that already exist elsewhere in the list. This matches how the
# source workbook was edited: cell A3 ("Construtora Beta") now holds the
# same value + formatting as A14 ("JAPJ CONSTRUCOES CIVIS LTDA PU_SUDESTE"),
# and cell A5 ("Construtora Delta") now holds the same value + formatting as
# A2 ("FG EMPREIT. MAO DE OBRA LTDA PU_SUL").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "JAPJ CONSTRUCOES CIVIS LTDA PU_SUDESTE"
$ws.Range("A3").WrapText = $false
$ws.Range("A3").VerticalAlignment = -4107

$ws.Range("A5").Value = "FG EMPREIT. MAO DE OBRA LTDA PU_SUL"
$ws.Range("A5").WrapText = $false
$ws.Range("A5").VerticalAlignment = -4107

$ws.Range("C3").Select() | Out-Null

$wb.Save()
